# Apply updated crypto price/volume figures to the active sheet.
# Values are written as literal text (matching the workbook's original
# inline-string cell type) by temporarily forcing a text number format
# before assignment, then clearing the format so the cell's style stays
# identical to before the edit (no stray format index is left behind).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" '42.972.93'
Set-TextValue "E2" '  -0.27%  '
Set-TextValue "D3" '2.220.64'
Set-TextValue "E3" '  -1.25%  '
Set-TextValue "E4" '  -0.17%  '
Set-TextValue "D5" '257.02'
Set-TextValue "E5" '  +4.64%  '
Set-TextValue "D6" '0.618'
Set-TextValue "E6" '  +0.47%  '
Set-TextValue "D7" '77.35'
Set-TextValue "E7" '  +2.63%  '
Set-TextValue "E8" '  -0.05%  '
Set-TextValue "D9" '0.598'
Set-TextValue "E9" '  -1.38%  '
Set-TextValue "D10" '42.47'
Set-TextValue "E10" '  +3.58%  '
Set-TextValue "D11" '0.0913'
Set-TextValue "E11" '  -2.18%  '
Set-TextValue "D12" '7.03'
Set-TextValue "E12" '  +0.94%  '
Set-TextValue "E13" '  +1.02%  '
Set-TextValue "D14" '2.553.42'
Set-TextValue "E14" '  -1.34%  '
Set-TextValue "D15" '14.54'
Set-TextValue "E15" '  -0.76%  '
Set-TextValue "D16" '2.217.01'
Set-TextValue "E16" '  -1.10%  '
Set-TextValue "D17" '0.786'
Set-TextValue "E17" '  -1.17%  '
Set-TextValue "D18" '42.904.56'
Set-TextValue "E18" '  -0.16%  '
Set-TextValue "E19" '  -1.34%  '
Set-TextValue "D20" '71.21'
Set-TextValue "E20" '  +0.07%  '
Set-TextValue "D21" '5.99'
Set-TextValue "E21" '  +0.18%  '
Set-TextValue "D22" '2.23'
Set-TextValue "E22" '  +1.59%  '
Set-TextValue "D23" '231.01'
Set-TextValue "E23" '  +0.26%  '
Set-TextValue "D24" '9.37'
Set-TextValue "E24" '  -5.44%  '
Set-TextValue "E25" '  -0.13%  '
Set-TextValue "D26" '43.06'
Set-TextValue "E26" '  +11.16%  '
Set-TextValue "D27" '10.82'
Set-TextValue "E27" '  -0.71%  '
Set-TextValue "E28" '  -2.57%  '
Set-TextValue "D29" '2.21'
Set-TextValue "E29" '  -1.93%  '
Set-TextValue "D30" '2.21'
Set-TextValue "E30" '  +3.65%  '
Set-TextValue "D31" '172.99'
Set-TextValue "E31" '  +0.16%  '
Set-TextValue "D32" '20.45'
Set-TextValue "E32" '  +0.65%  '
Set-TextValue "D33" '0.0878'
Set-TextValue "E33" '  +10.33%  '
Set-TextValue "D34" '5.25'
Set-TextValue "E34" '  -0.92%  '
Set-TextValue "E36" '  +8.57%  '
Set-TextValue "E37" '  -0.48%  '
Set-TextValue "D39" '12.97'
Set-TextValue "E39" '  -1.15%  '
Set-TextValue "B40" 'NEARProtocol'
Set-TextValue "C40" 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue "D40" '2.82'
Set-TextValue "E40" '  +18.02%  '
Set-TextValue "B41" 'LidoDAOToken'
Set-TextValue "C41" 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue "D41" '2.12'
Set-TextValue "E41" '  -0.46%  '
Set-TextValue "D42" '0.203'
Set-TextValue "E42" '  -1.69%  '
Set-TextValue "D43" '5.31'
Set-TextValue "E43" '  -3.32%  '
Set-TextValue "D44" '60.36'
Set-TextValue "E44" '  +1.22%  '
Set-TextValue "D45" '103.21'
Set-TextValue "E45" '  -2.09%  '
Set-TextValue "D46" '8.40'
Set-TextValue "E46" '  -3.76%  '
Set-TextValue "B47" 'Cronos'
Set-TextValue "C47" 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue "D47" '0.0979'
Set-TextValue "E47" '  -1.41%  '
Set-TextValue "B48" 'WOONetwork'
Set-TextValue "C48" 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
Set-TextValue "D48" '0.465'
Set-TextValue "E48" '  -3.02%  '
Set-TextValue "E49" '  +0.68%  '
Set-TextValue "E50" '  -0.79%  '
Set-TextValue "D51" '2.441.87'
Set-TextValue "E51" '  -0.76%  '
